# Applies the Dec 2 2024 cryptos-list refresh: updated prices/volumes for
# existing rows, plus a re-ranking swap among rows 43-46 (Kaspa/Algorand/
# EnergySwap/VeChain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while keeping the cell as TEXT (these Price/Volume
# columns are plain strings in the source data, e.g. "1.00" or "0.0492" -
# left alone, Excel would silently reinterpret them as numbers and mangle
# the formatting / introduce float rounding noise).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '95.963.85'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '3.665.03'
$ws.Range('E3').Value = '  -1.19%  '
Set-TextValue 'D4' '2.41'
$ws.Range('E4').Value = '  +28.13%  '
Set-TextValue 'D6' '228.12'
$ws.Range('E6').Value = '  -3.40%  '
Set-TextValue 'D7' '644.35'
$ws.Range('E7').Value = '  -1.01%  '
Set-TextValue 'D8' '0.429'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  +7.31%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '3.666.41'
$ws.Range('E11').Value = '  -1.12%  '
Set-TextValue 'D12' '47.57'
$ws.Range('E12').Value = '  +7.42%  '
Set-TextValue 'D13' '0.209'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E14').Value = '  -3.21%  '
Set-TextValue 'D15' '6.63'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '4.362.83'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = '95.703.83'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '3.669.11'
$ws.Range('E19').Value = '  -1.03%  '
Set-TextValue 'D20' '19.59'
$ws.Range('E20').Value = '  +5.14%  '
Set-TextValue 'D21' '12.93'
$ws.Range('E21').Value = '  -1.40%  '
Set-TextValue 'D22' '0.532'
$ws.Range('E22').Value = '  +5.60%  '
Set-TextValue 'D23' '523.43'
$ws.Range('E23').Value = '  +0.61%  '
Set-TextValue 'D24' '3.29'
$ws.Range('E24').Value = '  -2.92%  '
Set-TextValue 'D25' '0.246'
$ws.Range('E25').Value = '  +36.55%  '
Set-TextValue 'D26' '121.68'
$ws.Range('E26').Value = '  +19.95%  '
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('D29').Value = '3.860.31'
$ws.Range('E29').Value = '  -1.18%  '
Set-TextValue 'D30' '12.97'
$ws.Range('E30').Value = '  -2.69%  '
Set-TextValue 'D31' '13.08'
$ws.Range('E31').Value = '  +7.79%  '
Set-TextValue 'D32' '2.98'
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -1.28%  '
Set-TextValue 'D35' '1.81'
$ws.Range('E35').Value = '  -3.04%  '
Set-TextValue 'D36' '32.73'
$ws.Range('E36').Value = '  +1.45%  '
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  +0.29%  '
Set-TextValue 'D38' '0.608'
$ws.Range('E38').Value = '  +3.35%  '
Set-TextValue 'D39' '612.59'
$ws.Range('E39').Value = '  -4.98%  '
Set-TextValue 'D41' '8.43'
$ws.Range('E41').Value = '  -4.08%  '
Set-TextValue 'D42' '7.02'
$ws.Range('E42').Value = '  +3.06%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D43' '0.162'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.491'
$ws.Range('E44').Value = '  +13.14%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '40.08'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0492'
$ws.Range('E46').Value = '  +8.63%  '
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('E48').Value = '  -1.21%  '
Set-TextValue 'D49' '8.96'
$ws.Range('E49').Value = '  +5.67%  '
$ws.Range('E50').Value = '  -0.50%  '
Set-TextValue 'D51' '23.50'
$ws.Range('E51').Value = '  -0.30%  '
